# Applies the "Added note about misting to bullet on spider mites" commit.
#
# Most hunks in the target diff are pure run-merges (adjacent <w:r> runs with
# identical formatting collapsed into one run, no visible text change) which
# is exactly what Word's Find/Replace does when a match spans runs. The two
# exceptions are:
#   - the "phages" bullet, where wording genuinely changes
#   - the "mites can ..." bullet, where a new sentence is appended
#   - the lastRenderedPageBreak bookmark moving from one run/paragraph to another
#
# wdFindContinue = 1 (w:Wrap), wdReplaceAll = 2

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Table cell describing Umbilicus rupestris: merge 3 runs -> 1 (no text change).
Replace-Text "Flat green circular fleshy cap on a pale stem" "Flat green circular fleshy cap on a pale stem"

# 2. "Bonsai Pests: Pathogens" heading: merge 2 runs -> 1 (no text change).
Replace-Text "Bonsai Pests: Pathogens" "Bonsai Pests: Pathogens"

# 3. Spider-mite bullet: merge "mites can" + "(a) hide..." runs, then append a
#    brand-new sentence about misting after the existing "... later."
Replace-Text "mites can (a) hide from the insecticide in it and (b) use it to reinfect the plant faster" "mites can (a) hide from the insecticide in it and (b) use it to reinfect the plant faster"

Replace-Text "use it to reinfect the plant faster later." "use it to reinfect the plant faster later.  Subsequent (non-insecticidal) misting of leaves can slow reinfection and highlight any left-over or new webbing."

# 4. "Phages" bullet: reworded.
Replace-Text "Some infestations are treatable by phages – e.g. nematodes that predate wine weevils – but these tend to be very specific (one nematode species per pest) so not really recommended for beginner bonsai practitioners." "Some infestations are treatable by phages, e.g. nematodes that predate wine weevils.  However, these tend to be very specific – one nematode species per pest – so not recommended for beginner bonsai practitioners."

# 5. "them: look for pinprick-sized spots." merge of split runs (no text change).
Replace-Text "them: look for pinprick-sized spots." "them: look for pinprick-sized spots."

# 6. Move <w:lastRenderedPageBreak/> from the "Pesticides: environmental
#    considerations" heading run onto the "Birds are very good..." run.
$rng = $d.Content
$rng.Find.Execute("Pesticides: e", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rng.Find.Found) {
    $rng.Characters.First.Delete()
    $rng.InsertBefore("Pesticides: e")
}

# 7. "can also kill beneficial insects and (b) are often produced..." merges.
Replace-Text "can also kill beneficial insects and " "can also kill beneficial insects and "
Replace-Text "are often produced by industrial processes with a high environmental footprint." "are often produced by industrial processes with a high environmental footprint."

# 8. "...where the leaves will tend to come off too." merge.
Replace-Text "where the leaves will tend to come off too." "where the leaves will tend to come off too."

# 9. "Thanks to environmental regulation, most modern household insecticides..." merge.
$rsquo = [char]0x2019
$regText = "Thanks to environmental regulation, most modern household insecticides will break down over time, preventing build-up in the ecosystem.  However, it" + $rsquo + "s good to check the "
Replace-Text $regText $regText

# 10. "as a "trap crop"" merge.
$ldq = [char]0x201C
$rdq = [char]0x201D
$trapText = "as a " + $ldq + "trap crop" + $rdq + " "
Replace-Text $trapText $trapText
